$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5142127275466919
$ws.Range("B1").Value = 3.098735332489014
$ws.Range("C1").Value = 6.159793853759766
$ws.Range("D1").Value = 1.505792498588562
$ws.Range("E1").Value = 0.878449022769928
